$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.104.10"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "2.264.81"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.97"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.23"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.12"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.74"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.35"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").Value = "2.607.61"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.847"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "2.266.37"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "42.979.17"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.80"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.99"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.48"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.59"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.58"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.83"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.24"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.12"
$ws.Range("E29").Value = "  -6.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.28"
$ws.Range("E31").Value = "  -4.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.43"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.22"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0899"
$ws.Range("E34").Value = "  -5.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.70"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.126"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("E40").Value = "  -7.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  +7.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.26"
$ws.Range("E42").Value = "  +8.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.81"
$ws.Range("E43").Value = "  +9.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.234"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.06"
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.36"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.58"
$ws.Range("E48").Value = "  -3.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0989"
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.24"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.28"
$ws.Range("E51").Value = "  -0.11%  "
